$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need the Text number format
# applied first, otherwise Excel auto-converts the assigned string into a
# numeric value (changing the cell from a text/string cell to a number cell).
$numericLookingCells = @(
    "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D14", "D16", "D17", "D19", "D21", "D22", "D23", "D24", "D25", "D28", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D40", "D41", "D42", "D44", "D46", "D47", "D48", "D49", "D50"
)
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Cell value updates (price / 1h volume refresh) ---
$ws.Range("D2").Value = "29.347.84"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "1.840.70"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("D4").Value = "0.9984"
$ws.Range("D5").Value = "239.98"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").Value = "0.6285"
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("D7").Value = "0.9996"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "0.07447"
$ws.Range("E8").Value = "  -1.81%  "
$ws.Range("D9").Value = "0.2900"
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("D10").Value = "24.81"
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("D11").Value = "0.07736"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "1.843.28"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("D14").Value = "0.6776"
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("E15").Value = "  -2.39%  "
$ws.Range("D16").Value = "81.96"
$ws.Range("E16").Value = "  -1.69%  "
$ws.Range("D17").Value = "6.233"
$ws.Range("D18").Value = "29.293.12"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("D19").Value = "229.14"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "7.419"
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("D23").Value = "0.9989"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").Value = "158.85"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").Value = "8.462"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("E26").Value = "  -3.26%  "
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("D28").Value = "0.06506"
$ws.Range("E28").Value = "  +14.43%  "
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").Value = "4.062"
$ws.Range("E31").Value = "  -2.24%  "
$ws.Range("D32").Value = "4.061"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").Value = "1.837"
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("D34").Value = "1.139"
$ws.Range("E34").Value = "  -1.61%  "
$ws.Range("D35").Value = "0.6942"
$ws.Range("E35").Value = "  -1.40%  "
$ws.Range("E36").Value = "  -0.67%  "
$ws.Range("D37").Value = "0.01854"
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("D38").Value = "2.812"
$ws.Range("E38").Value = "  +1.34%  "
$ws.Range("D39").Value = "1.238.13"
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("D40").Value = "6.773"
$ws.Range("E40").Value = "  +3.97%  "
$ws.Range("D41").Value = "0.9332"
$ws.Range("E41").Value = "  +2.50%  "
$ws.Range("D42").Value = "0.9993"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").Value = "1.992.77"
$ws.Range("E43").Value = "  -1.19%  "
$ws.Range("D44").Value = "100.72"
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.00000000118"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "7.050"
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.712"
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("D49").Value = "0.1152"
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("D50").Value = "8.998"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("E51").Value = "  -1.72%  "
